$d = $word.ActiveDocument

# --- 1. Resize the table's grid columns (tblGrid / tcW) ---
# New widths below are expressed in points (twips / 20), matching the
# target w:w values of 1524,1524,8693,1524,1529,1529,1530,1530,1530 (twips).
$t = $d.Tables.Item(1)
$widths = @(76.2, 76.2, 434.65, 76.2, 76.45, 76.45, 76.5, 76.5, 76.5)
for ($i = 1; $i -le $t.Columns.Count; $i++) {
    $t.Columns.Item($i).Width = $widths[$i - 1]
}

# --- 2. Update the three product URLs to the new search-query URLs ---
$d.Content.Find.Execute(
    "https://www.fairprice.com.sg/product/maggi-2-min-big-noodles-curry-5s-x-111g-13064199",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://www.fairprice.com.sg/search?query=Maggi%20Big%20Curry%20Noodle", 2)

$d.Content.Find.Execute(
    "https://www.fairprice.com.sg/product/lipton-yellow-label-teabags-100s-72207",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://www.fairprice.com.sg/search?query=Lipton%20Yellow%20Label%20Tea", 2)

$d.Content.Find.Execute(
    "https://www.fairprice.com.sg/product/pantene-shampoo-hair-fall-control-680ml-13194383",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://www.fairprice.com.sg/search?query=Panteen%20Shampoo%20Hail%20Fall%20Control", 2)
